# Inquiry_4.xlsx edit: replace Discount/Total Amount columns with a
# Status/Remark layout, refresh the product rows, and add a 5th row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: D1/E1 get relabeled, F1 (old "Status" header) goes away.
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Remark"
$ws.Range("F1").Clear()

# ---------------------------------------------------------------------
# 2. Fix up the D/E column styles before we overwrite values, so the
#    new text content picks up the right formatting (D becomes the
#    "text" style used by C, E (rows 3+) becomes the style used by A/B).
# ---------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("D2:D4").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("E3:E4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Overwrite the product rows with the new data.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "MUTTON 1KG"
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 2300
$ws.Range("D2").Value = "Available"
$ws.Range("E2").Value = "We can supply this all items tommorrow morning"

$ws.Range("A3").Value = "BAIRAHA CHICKEN SUSAGES (500g)"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 750
$ws.Range("D3").Value = "Available"
$ws.Range("E3").Value = ""

$ws.Range("A4").Value = "HALAL CHICKEN"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 850
$ws.Range("D4").Value = "Available"
$ws.Range("E4").Value = ""

# Old "Status" column (F) is no longer used.
$ws.Range("F2:F4").Clear()

# ---------------------------------------------------------------------
# 4. Add the new 5th row, copying formatting from row 4 (which now has
#    the correct final style pattern) and then filling in its values.
# ---------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A5").Value = "BEEF 1KG"
$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 2000
$ws.Range("D5").Value = "Available"
$ws.Range("E5").Value = ""

# ---------------------------------------------------------------------
# 5. Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 32.43
$ws.Columns.Item(3).ColumnWidth = 14.14
$ws.Columns.Item(4).ColumnWidth = 20.43
$ws.Columns.Item(5).ColumnWidth = 45.86
$ws.Columns.Item(6).ColumnWidth = 8.71
